# SP sync upsert: 2025-11-12T16:40:45.5314750Z - MSCA_DF_00 - Course List.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Amendment")
$ws.Activate()

# Row 5 (MSCA_DF_10 Introduction to AI for Financial Applications): proposed new institute
# resolved from "???" to "POZ", and the type of change is now marked as "Institute change"
$ws.Range("F5").Value = "POZ"
$ws.Range("G5").Value = "Institute change"

# Row 11 (MSCA_DF_33 Gender and Diversity Dimension in Research): previously marked to be
# dropped, now changed instead to an institute change (new EC 2, move to BBU)
$ws.Range("C11").Value = 2
$ws.Range("F11").Value = "BBU"
$ws.Range("G11").Value = "Institute change"

# Row 15 (MSCA_DF_42 Intellectual Property Rights and Patenting): previously marked to be
# dropped, now changed instead to an institute change (new EC 1, move to POZ)
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = "POZ"
$ws.Range("G15").Value = "Institute change"

# Row 19 (new course "Practice of Digital Finance"): course column set to N/A, matching row 18
$ws.Range("A19").Value = "N/A"

# Leave the Amendment sheet's selection on the last-edited cell, then return focus to master
$ws.Range("G11").Select()
$wb.Worksheets.Item("master").Activate()
